$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16's A cell used to carry the "last row" date style (plain date).
# That style now moves to the new last row (17), so A16 gets the regular
# date+time style used by the rest of column A.
$ws.Range("A16").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new record as row 17.
$ws.Range("A17").Value = 45856
$ws.Range("A17").NumberFormat = "YYYY-MM-DD"

$ws.Range("B17").Value = "haha"
$ws.Range("C17").Value = "haha"
$ws.Range("D17").Value = "2025-07-18 14:47:29"
$ws.Range("E17").Value = "2025-07-18 14:47:31"
$ws.Range("F17").Value = "2025-07-18 14:47:31"
$ws.Range("G17").Value = "2025-07-18 14:47:32"
$ws.Range("H17").Value = "2025-07-18 14:47:33"
$ws.Range("I17").Value = "2025-07-18 14:47:34"
$ws.Range("J17").Value = "2025-07-18 14:47:35"
$ws.Range("K17").Value = "0:00:01"
$ws.Range("L17").Value = "0:00:02"
$ws.Range("M17").Value = "0:00:06"
$ws.Range("N17").Value = "2025-07-18 14:47:37"
$ws.Range("O17").Value = "2025-07-18 14:47:38"
$ws.Range("P17").Value = "2025-07-18 14:47:38"
$ws.Range("Q17").Value = "2025-07-18 14:47:40"
$ws.Range("R17").Value = "2025-07-18 14:47:41"
$ws.Range("S17").Value = "0:00:02"
$ws.Range("T17").Value = "0:00:01"
$ws.Range("U17").Value = "0:00:04"
$ws.Range("V17").Value = "0:00:02"
